# Auto-generated Excel COM-interop edit script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")


# ---- ALC ----
$ws_ALC.Range("H21").Value = 66618.60000000001
$ws_ALC.Range("I21").Value = 66618.60000000001
$ws_ALC.Range("K21").Value = 66618.60000000001
$ws_ALC.Range("M21").Value = -66150.60000000001
$ws_ALC.Range("H23").Value = 66618.60000000001
$ws_ALC.Range("I23").Value = 66618.60000000001
$ws_ALC.Range("K23").Value = 66618.60000000001
$ws_ALC.Range("M23").Value = -66384.60000000001
$ws_ALC.Range("H62").Value = 8093.25
$ws_ALC.Range("I62").Value = 3602.2942
$ws_ALC.Range("K62").Value = 3602.2942
$ws_ALC.Range("M62").Value = -2978.2942
$ws_ALC.Range("H65").Value = 8093.25
$ws_ALC.Range("I65").Value = 3602.2942
$ws_ALC.Range("K65").Value = 18011.471
$ws_ALC.Range("M65").Value = -14891.471
$ws_ALC.Range("H96").Value = 2056.2
$ws_ALC.Range("J96").Value = 2791.7
$ws_ALC.Range("L96").Value = 8375.099999999999
$ws_ALC.Range("N96").Value = -11121.1
$ws_ALC.Range("H116").Value = 5893.385
$ws_ALC.Range("J116").Value = 6200.4
$ws_ALC.Range("L116").Value = 6200.4
$ws_ALC.Range("N116").Value = -13084.4
$ws_ALC.Range("H132").Value = 1465.4231
$ws_ALC.Range("I132").Value = 1109.1818
$ws_ALC.Range("K132").Value = 3327.5454
$ws_ALC.Range("M132").Value = -797.5454
$ws_ALC.Range("H138").Value = 2695.897
$ws_ALC.Range("J138").Value = 2897.6182
$ws_ALC.Range("L138").Value = 8692.854599999999
$ws_ALC.Range("N138").Value = -18972.8546

# ---- ARM ----
$ws_ARM.Range("H5").Value = 452.0909
$ws_ARM.Range("I5").Value = 173.66667
$ws_ARM.Range("J5").Value = 556.5
$ws_ARM.Range("K5").Value = 173.66667
$ws_ARM.Range("L5").Value = 556.5
$ws_ARM.Range("M5").Value = -61.66667000000001
$ws_ARM.Range("N5").Value = -780.5
$ws_ARM.Range("H32").Value = 5750914.5
$ws_ARM.Range("I32").Value = 6946108
$ws_ARM.Range("K32").Value = 6946108
$ws_ARM.Range("M32").Value = -6945821
$ws_ARM.Range("H61").Value = 28853898
$ws_ARM.Range("I61").Value = 31256590
$ws_ARM.Range("K61").Value = 31256590
$ws_ARM.Range("M61").Value = -31256378
$ws_ARM.Range("H88").Value = 4606.4585
$ws_ARM.Range("J88").Value = 4764.5264
$ws_ARM.Range("L88").Value = 4764.5264
$ws_ARM.Range("N88").Value = -5576.5264
$ws_ARM.Range("H91").Value = 4606.4585
$ws_ARM.Range("J91").Value = 4764.5264
$ws_ARM.Range("L91").Value = 4764.5264
$ws_ARM.Range("N91").Value = -7572.5264
$ws_ARM.Range("H92").Value = 70749.75
$ws_ARM.Range("J92").Value = 70749.75
$ws_ARM.Range("L92").Value = 70749.75
$ws_ARM.Range("N92").Value = -75741.75
$ws_ARM.Range("H132").Value = 4506.1763
$ws_ARM.Range("I132").Value = 4483.8335
$ws_ARM.Range("J132").Value = 4559.8
$ws_ARM.Range("K132").Value = 13451.5005
$ws_ARM.Range("L132").Value = 13679.4
$ws_ARM.Range("M132").Value = -10921.5005
$ws_ARM.Range("N132").Value = -18739.4
$ws_ARM.Range("H136").Value = 28853898
$ws_ARM.Range("I136").Value = 31256590
$ws_ARM.Range("K136").Value = 93769770
$ws_ARM.Range("M136").Value = -93767220

# ---- BSM ----
$ws_BSM.Range("H4").Value = 452.0909
$ws_BSM.Range("I4").Value = 173.66667
$ws_BSM.Range("J4").Value = 556.5
$ws_BSM.Range("K4").Value = 173.66667
$ws_BSM.Range("L4").Value = 556.5
$ws_BSM.Range("M4").Value = -58.66667000000001
$ws_BSM.Range("N4").Value = -786.5
$ws_BSM.Range("H86").Value = 4181.8184
$ws_BSM.Range("J86").Value = 4299.6665
$ws_BSM.Range("L86").Value = 4299.6665
$ws_BSM.Range("N86").Value = -6545.6665
$ws_BSM.Range("H89").Value = 4181.8184
$ws_BSM.Range("J89").Value = 4299.6665
$ws_BSM.Range("L89").Value = 21498.3325
$ws_BSM.Range("N89").Value = -32730.3325
$ws_BSM.Range("H107").Value = 2291.28
$ws_BSM.Range("I107").Value = 1920.1364
$ws_BSM.Range("K107").Value = 1920.1364
$ws_BSM.Range("M107").Value = -0.1364000000000942
$ws_BSM.Range("H132").Value = 89000
$ws_BSM.Range("J132").Value = 89000
$ws_BSM.Range("L132").Value = 89000
$ws_BSM.Range("N132").Value = -99120

# ---- CRP ----
$ws_CRP.Range("H16").Value = 987.55554
$ws_CRP.Range("I16").Value = 987.55554
$ws_CRP.Range("K16").Value = 987.55554
$ws_CRP.Range("M16").Value = -700.55554
$ws_CRP.Range("H31").Value = 637214.9
$ws_CRP.Range("I31").Value = 12514.267
$ws_CRP.Range("K31").Value = 12514.267
$ws_CRP.Range("M31").Value = -12219.267
$ws_CRP.Range("H34").Value = 637214.9
$ws_CRP.Range("I34").Value = 12514.267
$ws_CRP.Range("K34").Value = 12514.267
$ws_CRP.Range("M34").Value = -12312.267
$ws_CRP.Range("H99").Value = 2397.2307
$ws_CRP.Range("I99").Value = 2633.2
$ws_CRP.Range("J99").Value = 1610.6666
$ws_CRP.Range("K99").Value = 2633.2
$ws_CRP.Range("L99").Value = 1610.6666
$ws_CRP.Range("M99").Value = -1135.2
$ws_CRP.Range("N99").Value = -4606.6666
$ws_CRP.Range("H107").Value = 6706.5713
$ws_CRP.Range("I107").Value = 6491
$ws_CRP.Range("K107").Value = 6491
$ws_CRP.Range("M107").Value = -4571
$ws_CRP.Range("H113").Value = 987.55554
$ws_CRP.Range("I113").Value = 987.55554
$ws_CRP.Range("K113").Value = 987.55554
$ws_CRP.Range("M113").Value = 1182.44446
$ws_CRP.Range("H126").Value = 2397.2307
$ws_CRP.Range("I126").Value = 2633.2
$ws_CRP.Range("J126").Value = 1610.6666
$ws_CRP.Range("K126").Value = 7899.599999999999
$ws_CRP.Range("L126").Value = 4831.9998
$ws_CRP.Range("M126").Value = -5429.599999999999
$ws_CRP.Range("N126").Value = -9771.9998

# ---- CUL ----
$ws_CUL.Range("H37").Value = 61495
$ws_CUL.Range("J37").Value = 61495
$ws_CUL.Range("L37").Value = 184485
$ws_CUL.Range("N37").Value = -184709
$ws_CUL.Range("H61").Value = 1000153.4
$ws_CUL.Range("I61").Value = 3333398
$ws_CUL.Range("K61").Value = 10000194
$ws_CUL.Range("M61").Value = -9999979
$ws_CUL.Range("H95").Value = 5400
$ws_CUL.Range("I95").Value = 800
$ws_CUL.Range("K95").Value = 2400
$ws_CUL.Range("M95").Value = -341
$ws_CUL.Range("H117").Value = 1358.8
$ws_CUL.Range("H126").Value = 9666.5
$ws_CUL.Range("I126").Value = 8999.5
$ws_CUL.Range("K126").Value = 26998.5
$ws_CUL.Range("M126").Value = -22058.5
$ws_CUL.Range("H134").Value = 12118.516
$ws_CUL.Range("J134").Value = 12789.871
$ws_CUL.Range("L134").Value = 38369.613
$ws_CUL.Range("N134").Value = -48509.613
$ws_CUL.Range("H139").Value = 3945.9614
$ws_CUL.Range("J139").Value = 4832.778
$ws_CUL.Range("L139").Value = 14498.334
$ws_CUL.Range("N139").Value = -24778.334
$ws_CUL.Range("H140").Value = 179329.6
$ws_CUL.Range("I140").Value = 179329.6
$ws_CUL.Range("K140").Value = 537988.8
$ws_CUL.Range("M140").Value = -532808.8

# ---- GSM ----
$ws_GSM.Range("H64").Value = 99997.5
$ws_GSM.Range("J64").Value = 99997.5
$ws_GSM.Range("L64").Value = 99997.5
$ws_GSM.Range("N64").Value = -100493.5
$ws_GSM.Range("H67").Value = 99997.5
$ws_GSM.Range("J67").Value = 99997.5
$ws_GSM.Range("L67").Value = 99997.5
$ws_GSM.Range("N67").Value = -101713.5
$ws_GSM.Range("H75").Value = 49665.832
$ws_GSM.Range("J75").Value = 49665.832
$ws_GSM.Range("L75").Value = 49665.832
$ws_GSM.Range("N75").Value = -51413.832
$ws_GSM.Range("H78").Value = 49665.832
$ws_GSM.Range("J78").Value = 49665.832
$ws_GSM.Range("L78").Value = 148997.496
$ws_GSM.Range("N78").Value = -157733.496
$ws_GSM.Range("H92").Value = 60000
$ws_GSM.Range("J92").Value = 60000
$ws_GSM.Range("L92").Value = 60000
$ws_GSM.Range("N92").Value = -63744
$ws_GSM.Range("H113").Value = 3621.68
$ws_GSM.Range("I113").Value = 2371
$ws_GSM.Range("J113").Value = 4776.154
$ws_GSM.Range("K113").Value = 2371
$ws_GSM.Range("L113").Value = 4776.154
$ws_GSM.Range("M113").Value = -201
$ws_GSM.Range("N113").Value = -9116.154
$ws_GSM.Range("H132").Value = 111125224
$ws_GSM.Range("I132").Value = 125003380
$ws_GSM.Range("K132").Value = 375010140
$ws_GSM.Range("M132").Value = -375007610

# ---- LTW ----
$ws_LTW.Range("H13").Value = 28333
$ws_LTW.Range("J13").Value = 0
$ws_LTW.Range("L13").Value = 0
$ws_LTW.Range("N13").ClearContents()
$ws_LTW.Range("H98").Value = 90352.5
$ws_LTW.Range("J98").Value = 90352.5
$ws_LTW.Range("L98").Value = 90352.5
$ws_LTW.Range("N98").Value = -96342.5
$ws_LTW.Range("H133").Value = 68566.336
$ws_LTW.Range("J133").Value = 68566.336
$ws_LTW.Range("L133").Value = 68566.336
$ws_LTW.Range("N133").Value = -73626.336

# ---- WVR ----
$ws_WVR.Range("H12").Value = 4000000
$ws_WVR.Range("I12").Value = 4000000
$ws_WVR.Range("K12").Value = 4000000
$ws_WVR.Range("M12").Value = -3999858
$ws_WVR.Range("H98").Value = 93333
$ws_WVR.Range("J98").Value = 93333
$ws_WVR.Range("L98").Value = 93333
$ws_WVR.Range("N98").Value = -99323
$ws_WVR.Range("H107").Value = 23810106
$ws_WVR.Range("I107").Value = 29412298
$ws_WVR.Range("K107").Value = 88236894
$ws_WVR.Range("M107").Value = -88234974
$ws_WVR.Range("H129").Value = 99742.5
$ws_WVR.Range("J129").Value = 99742.5
$ws_WVR.Range("L129").Value = 99742.5
$ws_WVR.Range("N129").Value = -109742.5
$ws_WVR.Range("H132").Value = 402577.6
$ws_WVR.Range("I132").Value = 2783.4546
$ws_WVR.Range("J132").Value = 3334401.2
$ws_WVR.Range("K132").Value = 8350.363799999999
$ws_WVR.Range("L132").Value = 10003203.6
$ws_WVR.Range("M132").Value = -5820.363799999999
$ws_WVR.Range("N132").Value = -10008263.6
